$wb = $excel.ActiveWorkbook

# --- Reorder sheet tabs: "openings" moves to sit before "users" ---
# (before: sites, users, openings, groups -> after: sites, openings, users, groups)
$openings = $wb.Worksheets.Item("openings")
$users = $wb.Worksheets.Item("users")
$openings.Move($users)

# --- "groups" sheet: add two new columns (Log, MJ) next to the existing header row ---
$groups = $wb.Worksheets.Item("groups")
$groups.Range("C1").Value = "Log"
$groups.Range("D1").Value = "MJ"
$groups.Range("A1").Copy()
$groups.Range("C1:D1").PasteSpecial(-4122)   # xlPasteFormats - match existing header style
$null = $groups.Select()
$null = $groups.Range("F4").Select()

# --- re-select "users" as the active sheet/cell (re-fetch: index changed after Move) ---
$users = $wb.Worksheets.Item("users")
$null = $users.Select()
$null = $users.Range("A4").Select()
